$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row labels to append (A11:A20), matching the new shared strings added
$labels = @(
    "AMU2a",
    "AMU17a",
    "AMU16a",
    "AMU4a",
    "AMU18a",
    "AMU29a",
    "AMU45a",
    "AMU30a",
    "AMU28a",
    "AMU32a"
)

# Update two existing Factor values
$ws.Range("D3").Value = 1
$ws.Range("D9").Value = 0.1

# Append the new rows (11-20), each with the same B/C/D pattern as existing rows
$row = 11
foreach ($label in $labels) {
    $ws.Cells.Item($row, 1).Value = $label
    $ws.Cells.Item($row, 2).Value = 7
    $ws.Cells.Item($row, 3).Value = 8
    $ws.Cells.Item($row, 4).Value = 10
    $row++
}

# Update selection to match the saved view state
$ws.Range("F15").Select()
